$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2-17 (Sending cluster x Target cluster matrix, Fn1 -> Itga2)
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metric columns (see header row 1 for full labels).
$data = @(
    ,@("ECs", "Fn1", "Itga2", "ECs", 3, 1, 19.95578266666667, 59.867348, 0.0117373419656925, 0.0117373419656925, 3, 1, 1.740822, 5.222466, 0.4863878955914668, 0.4863878955914669, 34.739465493352, 312.655189440168, 0.005708901058530584, 0.005708901058530585)
    ,@("ECs", "Fn1", "Itga2", "FAPs", 3, 1, 19.95578266666667, 59.867348, 0.0117373419656925, 0.0117373419656925, 3, 1, 1.077748, 3.233244, 0.3011241710513264, 0.3011241710513265, 21.50730485743467, 193.565743716912, 0.0035343973697651, 0.0035343973697651)
    ,@("ECs", "Fn1", "Itga2", "M2", 3, 1, 19.95578266666667, 59.867348, 0.0117373419656925, 0.0117373419656925, 1, 0.3333333333333333, 0.03488166666666666, 0.104645, 0.009745982326006345, 0.009745982326006345, 0.6960909590511111, 6.26481863146, 0.0001143919273519317, 0.0001143919273519317)
    ,@("ECs", "Fn1", "Itga2", "sCs", 3, 1, 19.95578266666667, 59.867348, 0.0117373419656925, 0.0117373419656925, 3, 1, 0.7256300000000001, 2.17689, 0.2027419510312003, 0.2027419510312003, 14.48051457641334, 130.32463118772, 0.00237965161004488, 0.00237965161004488)
    ,@("FAPs", "Fn1", "Itga2", "ECs", 3, 1, 1637.343343333333, 4912.03003, 0.9630320723052701, 0.9630320723052702, 3, 1, 1.740822, 5.222466, 0.4863878955914668, 0.4863878955914669, 2850.32331362822, 25652.90982265398, 0.4684071430356497, 0.4684071430356498)
    ,@("FAPs", "Fn1", "Itga2", "FAPs", 3, 1, 1637.343343333333, 4912.03003, 0.9630320723052701, 0.9630320723052702, 3, 1, 1.077748, 3.233244, 0.3011241710513264, 0.3011241710513265, 1764.643513590813, 15881.79162231732, 0.2899922344687655, 0.2899922344687656)
    ,@("FAPs", "Fn1", "Itga2", "M2", 3, 1, 1637.343343333333, 4912.03003, 0.9630320723052701, 0.9630320723052702, 1, 0.3333333333333333, 0.03488166666666666, 0.104645, 0.009745982326006345, 0.009745982326006345, 57.11326472103888, 514.0193824893501, 0.009385693556064428, 0.009385693556064428)
    ,@("FAPs", "Fn1", "Itga2", "sCs", 3, 1, 1637.343343333333, 4912.03003, 0.9630320723052701, 0.9630320723052702, 3, 1, 0.7256300000000001, 2.17689, 0.2027419510312003, 0.2027419510312003, 1188.105450222967, 10692.9490520067, 0.1952470012447904, 0.1952470012447904)
    ,@("M2", "Fn1", "Itga2", "ECs", 3, 1, 17.50081933333334, 52.502458, 0.01029341242216722, 0.01029341242216722, 3, 1, 1.740822, 5.222466, 0.4863878955914668, 0.4863878955914669, 30.465811313492, 274.192301821428, 0.005006591206472976, 0.005006591206472978)
    ,@("M2", "Fn1", "Itga2", "FAPs", 3, 1, 17.50081933333334, 52.502458, 0.01029341242216722, 0.01029341242216722, 3, 1, 1.077748, 3.233244, 0.3011241710513264, 0.3011241710513265, 18.86147303486133, 169.753257313752, 0.00309959528291453, 0.003099595282914531)
    ,@("M2", "Fn1", "Itga2", "M2", 3, 1, 17.50081933333334, 52.502458, 0.01029341242216722, 0.01029341242216722, 1, 0.3333333333333333, 0.03488166666666666, 0.104645, 0.009745982326006345, 0.009745982326006345, 0.6104577463788889, 5.49411971741, 0.0001003194155407359, 0.0001003194155407359)
    ,@("M2", "Fn1", "Itga2", "sCs", 3, 1, 17.50081933333334, 52.502458, 0.01029341242216722, 0.01029341242216722, 3, 1, 0.7256300000000001, 2.17689, 0.2027419510312003, 0.2027419510312003, 12.69911953284667, 114.29207579562, 0.002086906517238975, 0.002086906517238975)
    ,@("sCs", "Fn1", "Itga2", "ECs", 3, 1, 25.39612333333333, 76.18836999999999, 0.01493717330687017, 0.01493717330687017, 3, 1, 1.740822, 5.222466, 0.4863878955914668, 0.4863878955914669, 44.21013021338, 397.8911719204199, 0.007265260290813613, 0.007265260290813614)
    ,@("sCs", "Fn1", "Itga2", "FAPs", 3, 1, 25.39612333333333, 76.18836999999999, 0.01493717330687017, 0.01493717330687017, 3, 1, 1.077748, 3.233244, 0.3011241710513264, 0.3011241710513265, 27.37062113025333, 246.33559017228, 0.00449794392988128, 0.00449794392988128)
    ,@("sCs", "Fn1", "Itga2", "M2", 3, 1, 25.39612333333333, 76.18836999999999, 0.01493717330687017, 0.01493717330687017, 1, 0.3333333333333333, 0.03488166666666666, 0.104645, 0.009745982326006345, 0.009745982326006345, 0.8858591087388888, 7.97273197865, 0.0001455774270492504, 0.0001455774270492504)
    ,@("sCs", "Fn1", "Itga2", "sCs", 3, 1, 25.39612333333333, 76.18836999999999, 0.01493717330687017, 0.01493717330687017, 3, 1, 0.7256300000000001, 2.17689, 0.2027419510312003, 0.2027419510312003, 18.42818897436667, 165.8537007693, 0.003028391659126024, 0.003028391659126024)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $item = $data[$i]
    for ($c = 0; $c -lt $item.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $item[$c]
    }
}

Write-Output "Updated rows 2..17 ($($data.Count) rows) across columns A..T"